{"js": "// Append the 08/08/22 stand-up notes entry to the end of the document body,\n// mirroring the pattern used for every previous week's entry.\nconst body = context.document.body;\n\nconst paragraphsToAdd = [\n  \"\",\n  \"\",\n  \"08/08/22\",\n  \"Devin\",\n  \"Yesterday/Friday worked on front end web pages and styling. Today plan on finishing some of the pages I started and looking into either bank end mongo, discussion board page, contact us email form, or the search function. Blockers - lack of knowledge.\",\n  \"\",\n  \"Hodan\",\n  \"Yesterday/Friday worked on content and styling on web pages and links on nav bar. Today will set up the search function and movies api. Blockers - lack of knowledge.\",\n  \"\",\n  \"Toseef\",\n  \"Yesterday/Friday worked on content and styling on some of the web pages. Today will be adding more content and styling for unfinished pages, working on setting up the back-end. Blockers - lack of knowledge.\",\n  \"\",\n  \"Waseem\",\n  \"Yesterday/Friday worked on content and styling. Today will be adding more content and working on styling. Blockers - \"\n];\n\nfor (const text of paragraphsToAdd) {\n  body.insertParagraph(text, Word.InsertLocation.end);\n}\n\nawait context.sync();\n", "ps1": "# Append the 08/08/22 stand-up notes entry to the end of the document,\n# mirroring the pattern used for every previous week's entry.\n$d = $word.ActiveDocument\n\n$paragraphsToAdd = @(\n  \"\",\n  \"\",\n  \"08/08/22\",\n  \"Devin\",\n  \"Yesterday/Friday worked on front end web pages and styling. Today plan on finishing some of the pages I started and looking into either bank end mongo, discussion board page, contact us email form, or the search function. Blockers - lack of knowledge.\",\n  \"\",\n  \"Hodan\",\n  \"Yesterday/Friday worked on content and styling on web pages and links on nav bar. Today will set up the search function and movies api. Blockers - lack of knowledge.\",\n  \"\",\n  \"Toseef\",\n  \"Yesterday/Friday worked on content and styling on some of the web pages. Today will be adding more content and styling for unfinished pages, working on setting up the back-end. Blockers - lack of knowledge.\",\n  \"\",\n  \"Waseem\",\n  \"Yesterday/Friday worked on content and styling. Today will be adding more content and working on styling. Blockers - \"\n)\n\nforeach ($text in $paragraphsToAdd) {\n  # Re-fetch the end-of-document range each time: ranges in this host do\n  # not live-track document growth, so reusing a stale range after a\n  # mutation concatenates text into the wrong paragraph.\n  $r = $d.Content\n  $r.Collapse(0)\n  $r.InsertParagraphAfter()\n\n  if ($text -ne \"\") {\n    $r2 = $d.Content\n    $r2.Collapse(0)\n    $r2.InsertAfter($text)\n  }\n}\n"}
